# edit.ps1 -- reproduce the commit:
#   1) the table on the "B1- TYPES OF FINANCIAL DOCUMENTS" slide gets a new
#      table style ({439AEEBE-...} -> {F82CF26F-...})
#   2) the deck's theme colour scheme (ppt/theme/theme2.xml, the theme used
#      by the slide master / all slide layouts / the presentation itself)
#      is switched from the custom "Integral / Red Violet" palette back to
#      the stock Office palette.
#
# NOTE on (2): PowerPoint's object model only exposes read/write access to
# the *slide*-side theme colours (Slide.ThemeColorScheme / Master.ColorScheme
# / NotesMaster.ColorScheme all resolve to the same underlying theme part,
# ppt/theme/theme2.xml). The notes-master-only theme part (theme1.xml) has
# no writable surface in the object model, so this script focuses on the
# reachable, visible part of the theme swap: the 12 scheme colours.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style swap
# ---------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq "{439AEEBE-784E-494C-B618-DEF2152B119A}") {
                $tbl.ApplyStyle("{F82CF26F-518F-442D-B35D-098A8AF84679}")
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme colour scheme swap (Integral/"Red Violet" -> stock "Office")
# ---------------------------------------------------------------------
# Order of ThemeColorScheme.Item(n): dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink. RGB is packed the standard VBA way (R | G<<8 | B<<16).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
